# Apply the "cibmtr-reporting-ig" update to the Metadata sheet (sheet1):
#  - Version 0.1.6 -> 0.1.7
#  - Status active -> draft
#  - Date refreshed to the new publication timestamp
#  - Contact split into the two real contact lines (org + individual)
#  - New "Jurisdiction" row inserted right after the Contact rows
#  - Description / Purpose / Copyright / Immutable rows shift down one row

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- simple value edits (rows 1-11 keep their row numbers) ---------------
$ws1.Cells.Item(3, 2).Value = "0.1.7"
$ws1.Cells.Item(6, 2).Value = "draft"
$ws1.Cells.Item(8, 2).Value = "2024-08-27T12:23:18-05:00"
$ws1.Cells.Item(10, 2).Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Cells.Item(11, 2).Value = "Bob Milius (bmilius@nmdp.org)"

# --- make room for the new "Jurisdiction" row at row 12 -------------------
# Shift the old rows 12-15 (Description, Purpose, Copyright, Immutable) down
# to 13-16, carrying both their values and their existing formatting, without
# fabricating any brand-new cell style.
for ($r = 15; $r -ge 12; $r--) {
    $srcRange = $ws1.Range("A" + $r + ":B" + $r)
    $dstRange = $ws1.Range("A" + ($r + 1) + ":B" + ($r + 1))
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
    $ws1.Cells.Item($r + 1, 1).Value = $ws1.Cells.Item($r, 1).Value()
    $ws1.Cells.Item($r + 1, 2).Value = $ws1.Cells.Item($r, 2).Value()
}

# --- populate the new row 12: Jurisdiction / (no value) -------------------
$ws1.Cells.Item(12, 1).Value = "Jurisdiction"
$ws1.Cells.Item(12, 2).Value = ""

$excel.CutCopyMode = $false
